$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The author noticed the residual value in B19 looked "funny" because it was
# computed from a formula (3.52 - 3.01 = 0.51). They replaced it with a
# plain literal value (0.71).
$ws.Range("B19").Value = 0.71

# Reflect where the user ended up clicking/selecting afterwards.
$ws.Range("B20").Select()
